$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "0.1uF"
$ws.Range("E5").Value = 13
$ws.Range("F5").Value = "0.1uF"
$ws.Range("K5").Value = "C1, C1, C1, C3, C3, C3, C5, C5, C5, C10, C12, C13, C14"
$ws.Range("A6").Value = "10uF"
$ws.Range("E6").Value = 9
$ws.Range("F6").Value = "10uF"
$ws.Range("K6").Value = "C2, C2, C2, C11, C16, C17, C18, C18, C18"
$ws.Range("A7").Value = "47pF"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = "47pF"
$ws.Range("K7").Value = "C4"
$ws.Range("A8").Value = "Cap"
$ws.Range("K8").Value = "C6, C7"
$ws.Range("A9").Value = "100nF"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = "100nF"
$ws.Range("K9").Value = "C8, C9"
$ws.Range("A18").Value = "330R"
$ws.Range("E18").Value = 7
$ws.Range("F18").Value = "330R"
$ws.Range("K18").Value = "R1, R1, R1, R6, R7, R9, R11"
$ws.Range("A20").Value = "10KR"
$ws.Range("F20").Value = "10KR"
$ws.Range("K20").Value = "R3"
$ws.Range("A21").Value = "1KR"
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = "1KR"
$ws.Range("K21").Value = "R4"
$ws.Range("A22").Value = "360KR"
$ws.Range("E22").Value = 3
$ws.Range("F22").Value = "360KR"
$ws.Range("K22").Value = "R5"
$ws.Range("A24").Value = "430KR"
$ws.Range("F24").Value = "430KR"
